$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 150 (shifts old rows 150..236 down to 151..237)
$ws.Rows.Item(150).Insert()

# Populate the newly inserted row 150 with the new record's data
$ws.Cells.Item(150, 1).Value = 2
$ws.Cells.Item(150, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(150, 3).Value = "Coquimbo"
$ws.Cells.Item(150, 4).Value = 44910
$ws.Cells.Item(150, 5).Value = 4
$ws.Cells.Item(150, 6).Value = "Fruta"
$ws.Cells.Item(150, 7).Value = 100109
$ws.Cells.Item(150, 8).Value = "Uva"
$ws.Cells.Item(150, 9).Value = 100109001
$ws.Cells.Item(150, 10).Value = "Uva"
$ws.Cells.Item(150, 11).Value = "Flame Seedless"
$ws.Cells.Item(150, 12).Value = "Primera"
$ws.Cells.Item(150, 13).Value = 240
$ws.Cells.Item(150, 14).Value = 14000
$ws.Cells.Item(150, 15).Value = 15000
$ws.Cells.Item(150, 16).Value = 14500
$ws.Cells.Item(150, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(150, 18).Value = "Provincia de Huasco"
$ws.Cells.Item(150, 19).Value = 1450
$ws.Cells.Item(150, 20).Value = 10
